$wb = $excel.ActiveWorkbook

# New values for each of the three result sheets (sigma_010, sigma_025, sigma_050).
# Row key = spreadsheet row number; tuple = (A, B, C)
$sheetData = @{
    "sigma_010" = @{
        2  = @(1, 27.74685210548187, 29.64072440206422)
        3  = @(2, 27.72601770985651, 29.64624443768461)
        4  = @(3, 27.7472719319219, 29.65172153444765)
        5  = @(4, 27.74181502452763, 29.68502840867663)
        6  = @(5, 27.76465161124668, 29.66746913951022)
        7  = @(6, 27.75525077954714, 29.69195915441597)
        8  = @(7, 27.72720596662388, 29.67984825163132)
        9  = @(8, 27.75242302606005, 29.68480882278876)
        10 = @(9, 27.70276239876468, 29.64011819060485)
        11 = @(10, 27.74171756356102, 29.66125357605023)
        12 = @(27.74059681175914, 29.66491759178745)
    }
    "sigma_025" = @{
        2  = @(1, 19.64757907530271, 25.75139514928111)
        3  = @(2, 19.6786220366336, 25.73981064852625)
        4  = @(3, 19.64580210764992, 25.73641512225743)
        5  = @(4, 19.66381294214744, 25.77400236892506)
        6  = @(5, 19.65688034522016, 25.70572481385213)
        7  = @(6, 19.66018159374165, 25.72416372758671)
        8  = @(7, 19.64581118505269, 25.76654029806438)
        9  = @(8, 19.65858708914106, 25.78250589836954)
        10 = @(9, 19.66137067708744, 25.77153838656623)
        11 = @(10, 19.65273667786553, 25.76661543197593)
        12 = @(19.65713837298422, 25.75187118454048)
    }
    "sigma_050" = @{
        2  = @(1, 14.47400272728167, 21.46367997847292)
        3  = @(2, 14.47863409531755, 21.45417981805793)
        4  = @(3, 14.45419916273521, 21.43588377918806)
        5  = @(4, 14.48091838115517, 21.4336331281387)
        6  = @(5, 14.4808618836024, 21.42633938928358)
        7  = @(6, 14.46120866134584, 21.42286026209405)
        8  = @(7, 14.50787575119969, 21.48007564709545)
        9  = @(8, 14.48784222368089, 21.46533610255921)
        10 = @(9, 14.48527876073434, 21.47840060112332)
        11 = @(10, 14.49820341914615, 21.48126964730731)
        12 = @(14.48090250661989, 21.45416583533205)
    }
}

foreach ($sheetName in $sheetData.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $sheetData[$sheetName]
    foreach ($r in $rows.Keys) {
        $vals = $rows[$r]
        if ($r -eq 12) {
            # Row 12 keeps its "Média" label in column A; only B/C change.
            $ws.Cells.Item($r, 2).Value = $vals[0]
            $ws.Cells.Item($r, 3).Value = $vals[1]
        } else {
            $ws.Cells.Item($r, 1).Value = $vals[0]
            $ws.Cells.Item($r, 2).Value = $vals[1]
            $ws.Cells.Item($r, 3).Value = $vals[2]
        }
    }
}
